$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A25").Value = "Carglass"
$ws.Range("B25").Value = "Carglass"
$ws.Range("C25").Value = "Factuur"

$ws.Range("A26").Value = "Billit"
$ws.Range("B26").Value = "BILLIT"
$ws.Range("C26").Value = "Factuur"

$ws.Range("B27").Select()
